# Add I0 and IF columns (I and J) to the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row - copy formatting from an existing header cell (H1) so the
# new header cells share the same style as the rest of the header row.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for columns I (I0) and J (IF), keyed by row number
$data = @{
    2  = @(1, 3)
    3  = @(8, 8)
    4  = @(1, 7)
    5  = @(1, 4)
    6  = @(1, 8)
    7  = @(1, 5)
    8  = @(1, 5)
    9  = @(1, 6)
    10 = @(1, 6)
    11 = @(1, 6)
    12 = @(1, 5)
    13 = @(1, 5)
    14 = @(1, 7)
    15 = @(1, 5)
    16 = @(1, 5)
    17 = @(1, 6)
    18 = @(1, 6)
    19 = @(1, 6)
    20 = @(1, 4)
    21 = @(1, 5)
    22 = @(1, 1)
    23 = @(6, 6)
    24 = @(4, 4)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
